$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the header row strings -----------------------------------
# "_old" columns describe the FV2404 release, "_new" columns describe the
# FV2410 release; only the header labels change (values are untouched).
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$headerRange = $ws.Range("A1:U1")

# Remember the formatting so it can be restored after the range becomes a
# table header (ListObjects.Add snapshots the *current* header formatting
# into a table-level style override the moment it runs, so strip it first
# and re-apply afterwards via HeaderRowRange to avoid introducing a
# needless dxf/style diff).
$headerRange.ClearFormats()

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Turn the data range into a native Excel table --------------------
$tableRange = $ws.Range("A1:U72")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1)
$listObject.Name = "Table1"

# Restore the header row's original look (bold, centered, wrapped, grey
# fill, thin border on all sides) through the table's own HeaderRowRange.
$hrr = $listObject.HeaderRowRange
$hrr.Font.Bold = $true
$hrr.Interior.Color = 14277081
$hrr.Borders.LineStyle = 1
$hrr.HorizontalAlignment = -4108
$hrr.WrapText = $true

# --- Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
